# fix report with retry issue. read retryLimit from excel.
#
# Adds a new "retryLimit" parameter (value 1) to the Parameters sheet so the
# automation framework can read a retry limit from the workbook instead of
# having it hard-coded.

$wb = $excel.ActiveWorkbook

$wsScenarios  = $wb.Worksheets.Item("Scenarios")
$wsParameters = $wb.Worksheets.Item("Parameters")

# Append the new parameter row right after the existing "saucelab_url" row.
$wsParameters.Range("A12").Value = "retryLimit"
$wsParameters.Range("B12").Value = "1"

# Restore the selections that end up on each sheet after making the edit
# (Parameters sheet selection moves to A14, Scenarios keeps focus/selection
# on A4:A12). Select the target sheet last so it stays the active tab.
$wsParameters.Range("A14").Select() | Out-Null
$wsScenarios.Range("A4:A12").Select() | Out-Null
